# Commit message: "deleted previous data off excel"
#
# The first data row (row 2) of the InvoiceLog table previously held a
# sample/previous invoice entry (date, from, bill-to, addresses, inv no,
# PO no). That data is removed, leaving the row blank (only the
# auto-numbering formula in column A remains), and the now-empty last
# table row (row 18) is removed entirely since it was only there to keep
# the table in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 held real invoice data in B:F and H:I - wipe it back to blank,
# matching the style already used by the rest of the (empty) rows in the
# table (copy format from C2, which already carries the plain "s=1"
# style, onto B2 so the one-off date-number-format cell goes away too).
$ws.Range("C2").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B2:I2").ClearContents()

# Match the selection that was left on the sheet after clearing the row
# (now spans one extra column, B2:I2, instead of B2:H2).
$ws.Range("B2:I2").Select() | Out-Null

# The table (and sheet) had an extra trailing blank row (row 18) that is
# no longer needed - delete it so the table shrinks back from A1:X18 to
# A1:X17.
$ws.Rows("18:18").Delete()
